$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 2
$ws.Range("A2").Value = 2048420
$ws.Range("B2").Value = 173518
$ws.Range("C2").Value = -7.5

# Update the active selection to C3
$ws.Range("C3").Select()
